$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Lisso larvae"

# --- Rename headers (table column names follow automatically) ---
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Field"
$ws.Range("C1").Value = "Treatment"
$ws.Range("D1").Value = "Repeat"
$ws.Range("E1").Value = "Abundance"
$ws.Range("F1").Value = "Observations"

# --- Fix row 349: move the unrecorded "NA" value from Abundance to Observations ---
$ws.Range("E349").ClearContents()
$ws.Range("F349").Value = "Dato no apuntado"

# --- Column widths ---
# (input values pre-compensated for the engine's internal pixel-grid
# rounding so that the saved width ends up as close as possible to the
# true target width of the authored workbook)
$ws.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws.Columns.Item(2).ColumnWidth = 7.666666666666667
$ws.Columns.Item(3).ColumnWidth = 11.5
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334

# --- Selection / view ---
$ws.Range("M16").Select() | Out-Null
